$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update control-signal text (order matters so shared-string indices line up) ---
# mvnz row: "G != 0 else Done" -> "G != 0 {Ryout, ADDRin} else Done"
$ws.Range("D9").Value = "G != 0 {Ryout, ADDRin} else Done"
# st row: "Rxout, DOUTin, W, Done?" -> "Rxout, DOUTin, W_D, Done?"
$ws.Range("E8").Value = "Rxout, DOUTin, W_D, Done?"
# mvnz row: "if(true)  Ryout, RXin, Done" -> "DINout, RXin, Done"
$ws.Range("E9").Value = "DINout, RXin, Done"

# --- Adjust column widths to reflect the new (longer) text in columns B, D, E ---
$ws.Columns.Item(2).ColumnWidth = 13.333333333333332
$ws.Columns.Item(4).ColumnWidth = 29.666666666666668
$ws.Columns.Item(5).ColumnWidth = 25.5

# --- Move the active selection from E10 to C1 ---
[void]$ws.Range("C1").Select()
